$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.131583
$ws.Range("H2").Value = 39.394749
$ws.Range("I2").Value = 0.005137566080660421
$ws.Range("J2").Value = 0.005152598994293116
$ws.Range("M2").Value = 279.6078796666666
$ws.Range("N2").Value = 838.823639
$ws.Range("O2").Value = 0.5726675140320879
$ws.Range("P2").Value = 0.5775260287976519
$ws.Range("Q2").Value = 3671.694079296845
$ws.Range("R2").Value = 33045.24671367161
$ws.Range("S2").Value = 0.002942117195587381
$ws.Range("T2").Value = 0.002975760035160879

$ws.Range("G3").Value = 13.131583
$ws.Range("H3").Value = 39.394749
$ws.Range("I3").Value = 0.005137566080660421
$ws.Range("J3").Value = 0.005152598994293116
$ws.Range("O3").Value = 0.001411539217074187
$ws.Range("P3").Value = 0.001423514724607417
$ws.Range("Q3").Value = 9.050173196547334
$ws.Range("R3").Value = 81.451558768926
$ws.Range("S3").Value = 0.0000072518760031623087281460633690688410979419131763279438018798828125
$ws.Range("T3").Value = 0.0000073348005383736199867404821628458222448898595757782459259033203125

$ws.Range("G4").Value = 13.131583
$ws.Range("H4").Value = 39.394749
$ws.Range("I4").Value = 0.005137566080660421
$ws.Range("J4").Value = 0.005152598994293116
$ws.Range("M4").Value = 55.12872433333333
$ws.Range("N4").Value = 165.386173
$ws.Range("O4").Value = 0.1129096560274583
$ws.Range("P4").Value = 0.1138675822543568
$ws.Range("Q4").Value = 723.9274192672864
$ws.Range("R4").Value = 6515.346773405578
$ws.Range("S4").Value = 0.000580080818985705
$ws.Range("T4").Value = 0.0005867139898063878

$ws.Range("G5").Value = 13.131583
$ws.Range("H5").Value = 39.394749
$ws.Range("I5").Value = 0.005137566080660421
$ws.Range("J5").Value = 0.005152598994293116
$ws.Range("M5").Value = 12.3225355
$ws.Range("N5").Value = 24.645071
$ws.Range("O5").Value = 0.02523790023288966
$ws.Range("P5").Value = 0.01696801249072354
$ws.Range("Q5").Value = 161.8143976886965
$ws.Range("R5").Value = 970.886386132179
$ws.Range("S5").Value = 0.0001296613801835857
$ws.Range("T5").Value = 0.00008742936409485515571182734362309929565526545047760009765625

$ws.Range("G6").Value = 13.131583
$ws.Range("H6").Value = 39.394749
$ws.Range("I6").Value = 0.005137566080660421
$ws.Range("J6").Value = 0.005152598994293116
$ws.Range("M6").Value = 140.5068483333333
$ws.Range("N6").Value = 421.520545
$ws.Range("O6").Value = 0.2877733904904901
$ws.Range("P6").Value = 0.2902148617326603
$ws.Range("Q6").Value = 1845.077340957578
$ws.Range("R6").Value = 16605.69606861821
$ws.Range("S6").Value = 0.001478454809900588
$ws.Range("T6").Value = 0.001495360804692621

$ws.Range("I7").Value = 0.9851579054959454
$ws.Range("J7").Value = 0.9880405533247757
$ws.Range("M7").Value = 279.6078796666666
$ws.Range("N7").Value = 838.823639
$ws.Range("O7").Value = 0.5726675140320879
$ws.Range("P7").Value = 0.5775260287976519
$ws.Range("Q7").Value = 704068.500918817
$ws.Range("R7").Value = 6336616.508269354
$ws.Range("S7").Value = 0.5641679286694217
$ws.Range("T7").Value = 0.5706191370526923

$ws.Range("I8").Value = 0.9851579054959454
$ws.Range("J8").Value = 0.9880405533247757
$ws.Range("O8").Value = 0.001411539217074187
$ws.Range("P8").Value = 0.001423514724607417
$ws.Range("S8").Value = 0.001390589018618193
$ws.Range("T8").Value = 0.001406490276167078

$ws.Range("I9").Value = 0.9851579054959454
$ws.Range("J9").Value = 0.9880405533247757
$ws.Range("M9").Value = 55.12872433333333
$ws.Range("N9").Value = 165.386173
$ws.Range("O9").Value = 0.1129096560274583
$ws.Range("P9").Value = 0.1138675822543568
$ws.Range("Q9").Value = 138817.2548828349
$ws.Range("R9").Value = 1249355.293945515
$ws.Range("S9").Value = 0.1112338402422784
$ws.Range("T9").Value = 0.1125057889763491

$ws.Range("I10").Value = 0.9851579054959454
$ws.Range("J10").Value = 0.9880405533247757
$ws.Range("M10").Value = 12.3225355
$ws.Range("N10").Value = 24.645071
$ws.Range("O10").Value = 0.02523790023288966
$ws.Range("P10").Value = 0.01696801249072354
$ws.Range("Q10").Value = 31028.84334785863
$ws.Range("R10").Value = 186173.0600871518
$ws.Range("S10").Value = 0.02486331693254921
$ws.Range("T10").Value = 0.01676508445015619

$ws.Range("I11").Value = 0.9851579054959454
$ws.Range("J11").Value = 0.9880405533247757
$ws.Range("M11").Value = 140.5068483333333
$ws.Range("N11").Value = 421.520545
$ws.Range("O11").Value = 0.2877733904904901
$ws.Range("P11").Value = 0.2902148617326603
$ws.Range("Q11").Value = 353804.2139327844
$ws.Range("R11").Value = 3184237.92539506
$ws.Range("S11").Value = 0.283502230633078
$ws.Range("T11").Value = 0.2867440525694109

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.21459
$ws.Range("H12").Value = 3.64377
$ws.Range("I12").Value = 0.0004751930049796235
$ws.Range("J12").Value = 0.0004765834562732061
$ws.Range("M12").Value = 279.6078796666666
$ws.Range("N12").Value = 838.823639
$ws.Range("O12").Value = 0.5726675140320879
$ws.Range("P12").Value = 0.5775260287976519
$ws.Range("Q12").Value = 339.6089345643367
$ws.Range("R12").Value = 3056.48041107903
$ws.Range("S12").Value = 0.0002721275968471186
$ws.Range("T12").Value = 0.0002752393508921241

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.21459
$ws.Range("H13").Value = 3.64377
$ws.Range("I13").Value = 0.0004751930049796235
$ws.Range("J13").Value = 0.0004765834562732061
$ws.Range("O13").Value = 0.001411539217074187
$ws.Range("P13").Value = 0.001423514724607417
$ws.Range("Q13").Value = 0.8370849015533335
$ws.Range("R13").Value = 7.533764113980001
$ws.Range("S13").Value = 0.0000006707535622080679241650206894032049120824012788943946361541748046875
$ws.Range("T13").Value = 0.0000006784235675092038931799820651591392817181258578784763813018798828125

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.21459
$ws.Range("H14").Value = 3.64377
$ws.Range("I14").Value = 0.0004751930049796235
$ws.Range("J14").Value = 0.0004765834562732061
$ws.Range("M14").Value = 55.12872433333333
$ws.Range("N14").Value = 165.386173
$ws.Range("O14").Value = 0.1129096560274583
$ws.Range("P14").Value = 0.1138675822543568
$ws.Range("Q14").Value = 66.95879728802333
$ws.Range("R14").Value = 602.62917559221
$ws.Range("S14").Value = 0.00005365387873890355152628472978904028423130512237548828125
$ws.Range("T14").Value = 0.0000542674059082549727008222550939109396495041437447071075439453125

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.21459
$ws.Range("H15").Value = 3.64377
$ws.Range("I15").Value = 0.0004751930049796235
$ws.Range("J15").Value = 0.0004765834562732061
$ws.Range("M15").Value = 12.3225355
$ws.Range("N15").Value = 24.645071
$ws.Range("O15").Value = 0.02523790023288966
$ws.Range("P15").Value = 0.01696801249072354
$ws.Range("Q15").Value = 14.966828392945
$ws.Range("R15").Value = 89.80097035767
$ws.Range("S15").Value = 0.000011992873651042780468157188111799626994979917071759700775146484375
$ws.Range("T15").Value = 0.00000808667403891595693130218969191247424532775767147541046142578125

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.21459
$ws.Range("H16").Value = 3.64377
$ws.Range("I16").Value = 0.0004751930049796235
$ws.Range("J16").Value = 0.0004765834562732061
$ws.Range("M16").Value = 140.5068483333333
$ws.Range("N16").Value = 421.520545
$ws.Range("O16").Value = 0.2877733904904901
$ws.Range("P16").Value = 0.2902148617326603
$ws.Range("Q16").Value = 170.6582129171833
$ws.Range("R16").Value = 1535.92391625465
$ws.Range("S16").Value = 0.0001367479021803506
$ws.Range("T16").Value = 0.0001383116018664018

$ws.Range("G17").Value = 22.3716355
$ws.Range("H17").Value = 44.743271
$ws.Range("I17").Value = 0.008752619978390915
$ws.Range("J17").Value = 0.005852153878578699
$ws.Range("M17").Value = 279.6078796666666
$ws.Range("N17").Value = 838.823639
$ws.Range("O17").Value = 0.5726675140320879
$ws.Range("P17").Value = 0.5775260287976519
$ws.Range("Q17").Value = 6255.285566830527
$ws.Range("R17").Value = 37531.71340098316
$ws.Range("S17").Value = 0.005012341124292712
$ws.Range("T17").Value = 0.003379771189408332

$ws.Range("G18").Value = 22.3716355
$ws.Range("H18").Value = 44.743271
$ws.Range("I18").Value = 0.008752619978390915
$ws.Range("J18").Value = 0.005852153878578699
$ws.Range("O18").Value = 0.001411539217074187
$ws.Range("P18").Value = 0.001423514724607417
$ws.Range("Q18").Value = 15.41833729909233
$ws.Range("R18").Value = 92.510023794554
$ws.Range("S18").Value = 0.000012354666351645800017873370879595285032337415032088756561279296875
$ws.Range("T18").Value = 0.000008330627216825184537594799205617590587280574254691600799560546875

$ws.Range("G19").Value = 22.3716355
$ws.Range("H19").Value = 44.743271
$ws.Range("I19").Value = 0.008752619978390915
$ws.Range("J19").Value = 0.005852153878578699
$ws.Range("M19").Value = 55.12872433333333
$ws.Range("N19").Value = 165.386173
$ws.Range("O19").Value = 0.1129096560274583
$ws.Range("P19").Value = 0.1138675822543568
$ws.Range("Q19").Value = 1233.319726365314
$ws.Range("R19").Value = 7399.918358191882
$ws.Range("S19").Value = 0.0009882553110991774
$ws.Range("T19").Value = 0.0006663706131342134

$ws.Range("G20").Value = 22.3716355
$ws.Range("H20").Value = 44.743271
$ws.Range("I20").Value = 0.008752619978390915
$ws.Range("J20").Value = 0.005852153878578699
$ws.Range("M20").Value = 12.3225355
$ws.Range("N20").Value = 24.645071
$ws.Range("O20").Value = 0.02523790023288966
$ws.Range("P20").Value = 0.01696801249072354
$ws.Range("Q20").Value = 275.6752726418102
$ws.Range("R20").Value = 1102.701090567241
$ws.Range("S20").Value = 0.0002208977497910268
$ws.Range("T20").Value = 0.000099299420109359598529165380487171432832838036119937896728515625

$ws.Range("G21").Value = 22.3716355
$ws.Range("H21").Value = 44.743271
$ws.Range("I21").Value = 0.008752619978390915
$ws.Range("J21").Value = 0.005852153878578699
$ws.Range("M21").Value = 140.5068483333333
$ws.Range("N21").Value = 421.520545
$ws.Range("O21").Value = 0.2877733904904901
$ws.Range("P21").Value = 0.2902148617326603
$ws.Range("Q21").Value = 3143.367996167116
$ws.Range("R21").Value = 18860.20797700269
$ws.Range("S21").Value = 0.002518771126856353
$ws.Range("T21").Value = 0.001698382028709969

$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1.218481333333333
$ws.Range("H22").Value = 3.655444
$ws.Range("I22").Value = 0.000476715440023584
$ws.Range("J22").Value = 0.0004781103460792403
$ws.Range("M22").Value = 279.6078796666666
$ws.Range("N22").Value = 838.823639
$ws.Range("O22").Value = 0.5726675140320879
$ws.Range("P22").Value = 0.5775260287976519
$ws.Range("Q22").Value = 340.6969820267462
$ws.Range("R22").Value = 3066.272838240716
$ws.Range("S22").Value = 0.0002729994459390188
$ws.Range("T22").Value = 0.0002761211694982146

$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1.218481333333333
$ws.Range("H23").Value = 3.655444
$ws.Range("I23").Value = 0.000476715440023584
$ws.Range("J23").Value = 0.0004781103460792403
$ws.Range("O23").Value = 0.001411539217074187
$ws.Range("P23").Value = 0.001423514724607417
$ws.Range("Q23").Value = 0.8397667747617779
$ws.Range("R23").Value = 7.557900972856
$ws.Range("S23").Value = 0.000000672902538978066236345491370196558733596248202957212924957275390625
$ws.Range("T23").Value = 0.0000006805971176309466009370498056074705317541884141974151134490966796875

$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1.218481333333333
$ws.Range("H24").Value = 3.655444
$ws.Range("I24").Value = 0.000476715440023584
$ws.Range("J24").Value = 0.0004781103460792403
$ws.Range("M24").Value = 55.12872433333333
$ws.Range("N24").Value = 165.386173
$ws.Range("O24").Value = 0.1129096560274583
$ws.Range("P24").Value = 0.1138675822543568
$ws.Range("Q24").Value = 67.17332153064578
$ws.Range("R24").Value = 604.5598937758119
$ws.Range("S24").Value = 0.0000538257763560412766834679942906660699009080417454242706298828125
$ws.Range("T24").Value = 0.0000544412691588369110503274617141045155221945606172084808349609375

$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1.218481333333333
$ws.Range("H25").Value = 3.655444
$ws.Range("I25").Value = 0.000476715440023584
$ws.Range("J25").Value = 0.0004781103460792403
$ws.Range("M25").Value = 12.3225355
$ws.Range("N25").Value = 24.645071
$ws.Range("O25").Value = 0.02523790023288966
$ws.Range("P25").Value = 0.01696801249072354
$ws.Range("Q25").Value = 15.01477948608733
$ws.Range("R25").Value = 90.088676916524
$ws.Range("S25").Value = 0.00001203129671479330993553828743625189190424862317740917205810546875
$ws.Range("T25").Value = 0.000008112582324216703944289551253543635311871184967458248138427734375

$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1.218481333333333
$ws.Range("H26").Value = 3.655444
$ws.Range("I26").Value = 0.000476715440023584
$ws.Range("J26").Value = 0.0004781103460792403
$ws.Range("M26").Value = 140.5068483333333
$ws.Range("N26").Value = 421.520545
$ws.Range("O26").Value = 0.2877733904904901
$ws.Range("P26").Value = 0.2902148617326603
$ws.Range("Q26").Value = 171.2049718996644
$ws.Range("R26").Value = 1540.84474709698
$ws.Range("S26").Value = 0.0001371860184747526
$ws.Range("T26").Value = 0.0001387547279803411

